$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dSF (column F) values, keyed by row number.
$values = @{
    2  = 1
    3  = -2
    5  = 6
    6  = -2
    7  = -2
    8  = -4
    10 = 1
    12 = 2
    13 = -1
    14 = -4
    16 = 1
    17 = 1
    18 = -4
    19 = 4
    20 = 5
    22 = -4
    23 = -1
    27 = 2
    28 = 3
    29 = -1
    30 = 5
    31 = 1
    32 = -3
    33 = -3
    34 = -3
    35 = 5
    36 = 4
    37 = 4
    38 = 2
    39 = -2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
